# Auto-generated update script for resum_diari_meteocat equivalent workbook
# Commit: Update automàtic: dades i banners [2026-02-19 18:20]
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$updates = @(
    @('E2', '2026-02-19 18:18:29'),
    @('H2', '68%'),
    @('I2', '2.3 mm'),
    @('E3', '2026-02-19 18:18:32'),
    @('I3', '3.8 mm'),
    @('E4', '2026-02-19 18:18:34'),
    @('H4', '58%'),
    @('J4', '1009.4 hPa'),
    @('L4', '39.2 km/h - 262º 17:59 TU'),
    @('O4', '11.7 °C'),
    @('E5', '2026-02-19 18:18:37'),
    @('I5', '6.8 mm'),
    @('E6', '2026-02-19 18:18:39'),
    @('H6', '72%'),
    @('J6', '1009.5 hPa'),
    @('O6', '10.5 °C'),
    @('E7', '2026-02-19 18:18:42'),
    @('J7', '1010.4 hPa'),
    @('E8', '2026-02-19 18:18:44'),
    @('J8', '1010.1 hPa'),
    @('E9', '2026-02-19 18:18:47'),
    @('E10', '2026-02-19 18:18:49'),
    @('O10', '11.0 °C'),
    @('E11', '2026-02-19 18:18:50'),
    @('H11', '68%'),
    @('O11', '5.3 °C'),
    @('E12', '2026-02-19 18:18:52'),
    @('O12', '11.1 °C'),
    @('E13', '2026-02-19 18:18:53'),
    @('J13', '1010.7 hPa'),
    @('E14', '2026-02-19 18:18:54'),
    @('E15', '2026-02-19 18:18:55'),
    @('E16', '2026-02-19 18:18:56'),
    @('E17', '2026-02-19 18:18:57'),
    @('E18', '2026-02-19 18:18:58'),
    @('J18', '1009.7 hPa'),
    @('L18', '26.6 km/h - 262º 17:50 TU'),
    @('O18', '11.8 °C'),
    @('E19', '2026-02-19 18:18:59'),
    @('E20', '2026-02-19 18:19:00'),
    @('H20', '88%'),
    @('E21', '2026-02-19 18:19:02'),
    @('J21', '1010.5 hPa'),
    @('O21', '6.4 °C'),
    @('E22', '2026-02-19 18:19:04'),
    @('E23', '2026-02-19 18:19:07'),
    @('I23', '7.8 mm'),
    @('E24', '2026-02-19 18:19:09'),
    @('H24', '64%'),
    @('J24', '1014.1 hPa'),
    @('E25', '2026-02-19 18:19:12'),
    @('I25', '3.6 mm'),
    @('E26', '2026-02-19 18:19:14'),
    @('H26', '58%'),
    @('K26', '8.9 MJ/m2'),
    @('O26', '3.0 °C'),
    @('E27', '2026-02-19 18:19:17'),
    @('E28', '2026-02-19 18:19:19'),
    @('H28', '68%'),
    @('J28', '1009.3 hPa'),
    @('E29', '2026-02-19 18:19:21'),
    @('H29', '74%'),
    @('E30', '2026-02-19 18:19:24'),
    @('J30', '1009.5 hPa'),
    @('E31', '2026-02-19 18:19:26'),
    @('J31', '1008.9 hPa'),
    @('E32', '2026-02-19 18:19:29'),
    @('H32', '67%'),
    @('E33', '2026-02-19 18:19:31'),
    @('J33', '1010.2 hPa'),
    @('L33', '51.1 km/h - 346º 17:38 TU'),
    @('E34', '2026-02-19 18:19:34'),
    @('L34', '76.7 km/h - 68º 17:59 TU'),
    @('E35', '2026-02-19 18:19:36'),
    @('J35', '1015.6 hPa'),
    @('E36', '2026-02-19 18:19:39'),
    @('J36', '1009.8 hPa'),
    @('E37', '2026-02-19 18:19:41'),
    @('H37', '72%'),
    @('J37', '1010.8 hPa'),
    @('L37', '49.3 km/h - 238º 17:50 TU'),
    @('O37', '5.7 °C'),
    @('E38', '2026-02-19 18:19:44'),
    @('O38', '11.9 °C'),
    @('E39', '2026-02-19 18:19:46'),
    @('I39', '4.3 mm'),
    @('E40', '2026-02-19 18:19:49'),
    @('H40', '75%'),
    @('J40', '1011.7 hPa'),
    @('O40', '6.2 °C'),
    @('E41', '2026-02-19 18:19:51'),
    @('J41', '1012.3 hPa'),
    @('E42', '2026-02-19 18:19:54'),
    @('H42', '74%'),
    @('E43', '2026-02-19 18:19:56'),
    @('E44', '2026-02-19 18:19:58'),
    @('I44', '7.9 mm'),
    @('E45', '2026-02-19 18:20:01'),
    @('I45', '3.3 mm'),
    @('J45', '1014.7 hPa'),
    @('O45', '2.5 °C'),
    @('E46', '2026-02-19 18:20:04'),
    @('J46', '1015.0 hPa'),
    @('O46', '12.9 °C'),
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = "'" + $u[1]
}
